$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2:F6").ClearContents()

$ws.Range("F2:F6").Select()
